$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D4"  = -7.769
    "D7"  = -7.585000000000001
    "E7"  = 16.925
    "E15" = 16.068
    "D16" = -8.635999999999999
    "E21" = 16.557
    "E22" = 16.468
    "E23" = 16.434
    "D28" = -8.055
    "D29" = -7.292
    "D32" = -7.935999999999998
    "E34" = 16.881
    "D40" = -8.132
    "E43" = 17.083
    "E45" = 16.774
    "E50" = 16.403
    "E51" = 16.532
    "D52" = -7.911000000000001
    "D57" = -8.101000000000001
    "D66" = -7.419
    "E66" = 17.325
    "E67" = 17.296
    "E79" = 17.07
    "E84" = 16.608
    "E92" = 17.829
    "E97" = 16.872
    "D100" = -8.180999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
